$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 151 & 152: Classification "Upcoming" -> "HWealth", clear StartTime ---
$ws.Range("E151").Value = "HWealth"
$ws.Range("I151").ClearContents()

$ws.Range("E152").Value = "HWealth"
$ws.Range("I152").ClearContents()

# --- Append new seminar rows 154-157 ---
# Seed formatting for the new rows by copying the matching cells from row 153
# (which already carries the correct date/time/text styles), then overwrite
# the copied values with the new seminar data.
foreach ($col in @("A","B","C","D","E","I")) {
    $src = $col + "153"
    $ws.Range($src).Copy($ws.Range($col + "154"))
    $ws.Range($src).Copy($ws.Range($col + "155"))
    $ws.Range($src).Copy($ws.Range($col + "156"))
    $ws.Range($src).Copy($ws.Range($col + "157"))
}

# Dates
$ws.Range("A154").Value = 45569
$ws.Range("A155").Value = 45583
$ws.Range("A156").Value = 45597
$ws.Range("A157").Value = 45611

# Speakers (column B) -- note: existing speaker "甘舒" reused for rows 154/155,
# existing speaker "陈焕" reused for row 156, and "Dave Vrane" is a brand-new name.
$ws.Range("B154").Value = "甘舒"
$ws.Range("B155").Value = "甘舒"
$ws.Range("B156").Value = "陈焕"
$ws.Range("B157").Value = "Dave Vrane"

# Titles (column C)
$ws.Range("C154").Value = "历史15: 以色列历史-1"
$ws.Range("C155").Value = "历史16:以色列历史-2"
$ws.Range("C156").Value = "退休后旅居生活， 从沿海到边陲，从繁华到中国最大无人区"
$ws.Range("C157").Value = "Recreational Auto Racing"

# Images (column D)
$ws.Range("D154").Value = "img/israel_his_1.jpg"
$ws.Range("D155").Value = "img/israel_history.jpg"
$ws.Range("D156").Value = "img/roadtrip.jpg"
$ws.Range("D157").Value = "img/auto_racing.jpg"

# Classification (column E) -- all new seminars are "Upcoming"
$ws.Range("E154").Value = "Upcoming"
$ws.Range("E155").Value = "Upcoming"
$ws.Range("E156").Value = "Upcoming"
$ws.Range("E157").Value = "Upcoming"

# StartTime (column I) -- 7:30 PM for all new seminars
$ws.Range("I154").Value = 0.8125
$ws.Range("I155").Value = 0.8125
$ws.Range("I156").Value = 0.8125
$ws.Range("I157").Value = 0.8125

# --- Restore selection to match the author's final cursor position ---
$null = $ws.Range("E152").Select()
